# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de
# handback has completed (in addition to the already-recorded zh-cn one),
# and that both locales are now "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdName   = "09ce23a3-7216-4962-86ab-66a9b48ad55e.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bdd6c27dde0f7bd76bb65608a94f3c41c13bc6a/e2e/09ce23a3-7216-4962-86ab-66a9b48ad55e.md"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: refresh the per-locale status column
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------
# zh-cn sheet: already-synced locale — refresh status + handback info
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, $null, $null, $mdName) | Out-Null
$wsZhCn.Range("J2").Value = "09ce23a3-7216-4962-86ab-66a9b48ad55e.4a6c2d377b311c634623ae3811ccf586e8ca18eb.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-10-17 16:40:28"
$wsZhCn.Columns.Item(3).ColumnWidth = 29.2
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------
# de-de sheet: newly-synced locale — refresh status + handback info
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, $null, $null, $mdName) | Out-Null
$wsDeDe.Range("J2").Value = "09ce23a3-7216-4962-86ab-66a9b48ad55e.4a6c2d377b311c634623ae3811ccf586e8ca18eb.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-10-17 16:41:06"
$wsDeDe.Columns.Item(3).ColumnWidth = 29.2
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
